$d = $word.ActiveDocument

# =====================================================================
# Stage 1: build the new "Knarot" block (13 paragraphs) as scratch
# content splitting off paragraph 2 ("Detta dokument..."), which is
# Normal-styled with NO explicit w:pStyle -- paragraphs split off of it
# inherit that same "no explicit style" state cleanly (unlike splitting
# off a Title/Heading paragraph, which would force an explicit pStyle
# on every paragraph we try to reset back to Normal).
# Italic formatting is applied in a SEPARATE second pass after all the
# text has been typed in, so that setting Font.Italic never "bleeds"
# into text inserted immediately afterwards (Word-style "sticky cursor
# formatting" on empty-range inserts).
# =====================================================================
$anchor = $d.Paragraphs(2)
$scratchStartIdx = 3
$anchor.Range.InsertParagraphAfter()
$italicRanges = @()

# --- paragraph 0  (style=Heading1) ---
$pIdx0 = $scratchStartIdx + 0
$pStart = $d.Paragraphs($pIdx0).Range.Start
$r = $d.Range($pStart, $pStart)
$r.InsertAfter('Knärot – ekologi samt krav på livsmiljön')
$d.Paragraphs($pIdx0).Range.InsertParagraphAfter()

# --- paragraph 1  (style=None) ---
$pIdx1 = $scratchStartIdx + 1
$pStart = $d.Paragraphs($pIdx1).Range.Start
$r = $d.Range($pStart, $pStart)
$r.InsertAfter('Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021).')
$d.Paragraphs($pIdx1).Range.InsertParagraphAfter()

# --- paragraph 2  (style=None) ---
$pIdx2 = $scratchStartIdx + 2
$pStart = $d.Paragraphs($pIdx2).Range.Start
$r = $d.Range($pStart, $pStart)
$r.InsertAfter('Samuel Johnsons doktorsavhandling ')
$r = $d.Range($r.End, $r.End)
$r.InsertAfter('“Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“')
$italicRanges += ,@($r.Start, $r.End)
$r = $d.Range($r.End, $r.End)
$r.InsertAfter(' (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: ')
$r = $d.Range($r.End, $r.End)
$r.InsertAfter('“Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” ')
$italicRanges += ,@($r.Start, $r.End)
$r = $d.Range($r.End, $r.End)
$r.InsertAfter('Vidare ')
$r = $d.Range($r.End, $r.End)
$r.InsertAfter('“More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”')
$italicRanges += ,@($r.Start, $r.End)
$d.Paragraphs($pIdx2).Range.InsertParagraphAfter()

# --- paragraph 3  (style=None) ---
$pIdx3 = $scratchStartIdx + 3
$pStart = $d.Paragraphs($pIdx3).Range.Start
$r = $d.Range($pStart, $pStart)
$r.InsertAfter('Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: ')
$r = $d.Range($r.End, $r.End)
$r.InsertAfter('“In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”')
$italicRanges += ,@($r.Start, $r.End)
$d.Paragraphs($pIdx3).Range.InsertParagraphAfter()

# --- paragraph 4  (style=None) ---
$pIdx4 = $scratchStartIdx + 4
$pStart = $d.Paragraphs($pIdx4).Range.Start
$r = $d.Range($pStart, $pStart)
$r.InsertAfter('En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022).')
$d.Paragraphs($pIdx4).Range.InsertParagraphAfter()

# --- paragraph 5  (style=None) ---
$pIdx5 = $scratchStartIdx + 5
$pStart = $d.Paragraphs($pIdx5).Range.Start
$r = $d.Range($pStart, $pStart)
$r.InsertAfter('Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022).')
$d.Paragraphs($pIdx5).Range.InsertParagraphAfter()

# --- paragraph 6  (style=Heading2) ---
$pIdx6 = $scratchStartIdx + 6
$pStart = $d.Paragraphs($pIdx6).Range.Start
$r = $d.Range($pStart, $pStart)
$r.InsertAfter('Referenser - knärot')
$d.Paragraphs($pIdx6).Range.InsertParagraphAfter()

# --- paragraph 7  (style=None) ---
$pIdx7 = $scratchStartIdx + 7
$pStart = $d.Paragraphs($pIdx7).Range.Start
$r = $d.Range($pStart, $pStart)
$r.InsertAfter('de Graaf M & Roberts M.R., 2009. ')
$r = $d.Range($r.End, $r.End)
$r.InsertAfter('Short-term response of the herbaceous layer within leave patches after harvest. ')
$italicRanges += ,@($r.Start, $r.End)
$r = $d.Range($r.End, $r.End)
$r.InsertAfter('Forest Ecology and Management 257, 1014-1025')
$d.Paragraphs($pIdx7).Range.InsertParagraphAfter()

# --- paragraph 8  (style=None) ---
$pIdx8 = $scratchStartIdx + 8
$pStart = $d.Paragraphs($pIdx8).Range.Start
$r = $d.Range($pStart, $pStart)
$r.InsertAfter('Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. ')
$r = $d.Range($r.End, $r.End)
$r.InsertAfter('Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. ')
$italicRanges += ,@($r.Start, $r.End)
$r = $d.Range($r.End, $r.End)
$r.InsertAfter('Ecological Applications, 22, 2049-2064 ')
$d.Paragraphs($pIdx8).Range.InsertParagraphAfter()

# --- paragraph 9  (style=None) ---
$pIdx9 = $scratchStartIdx + 9
$pStart = $d.Paragraphs($pIdx9).Range.Start
$r = $d.Range($pStart, $pStart)
$r.InsertAfter('Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. ')
$r = $d.Range($r.End, $r.End)
$r.InsertAfter('Interactive effects of drought and edge exposure on old-growth forest understory species. ')
$italicRanges += ,@($r.Start, $r.End)
$r = $d.Range($r.End, $r.End)
$r.InsertAfter('Landscape Ecology, 37, sid 1839-1853')
$d.Paragraphs($pIdx9).Range.InsertParagraphAfter()

# --- paragraph 10  (style=None) ---
$pIdx10 = $scratchStartIdx + 10
$pStart = $d.Paragraphs($pIdx10).Range.Start
$r = $d.Range($pStart, $pStart)
$r.InsertAfter('Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. ')
$r = $d.Range($r.End, $r.End)
$r.InsertAfter('Biological legacies buffer local species extinction after logging. ')
$italicRanges += ,@($r.Start, $r.End)
$r = $d.Range($r.End, $r.End)
$r.InsertAfter('Journal of Applied Ecology. 51, 53-62.')
$d.Paragraphs($pIdx10).Range.InsertParagraphAfter()

# --- paragraph 11  (style=None) ---
$pIdx11 = $scratchStartIdx + 11
$pStart = $d.Paragraphs($pIdx11).Range.Start
$r = $d.Range($pStart, $pStart)
$r.InsertAfter('Skogsstyrelsen, 2022. ')
$r = $d.Range($r.End, $r.End)
$r.InsertAfter('Vägledning för hänsyn till knärot. ')
$italicRanges += ,@($r.Start, $r.End)
$r = $d.Range($r.End, $r.End)
$r.InsertAfter('https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/')
$d.Paragraphs($pIdx11).Range.InsertParagraphAfter()

# --- paragraph 12  (style=None) ---
$pIdx12 = $scratchStartIdx + 12
$pStart = $d.Paragraphs($pIdx12).Range.Start
$r = $d.Range($pStart, $pStart)
$r.InsertAfter('SLU Artdatabanken, 2021. ')
$r = $d.Range($r.End, $r.End)
$r.InsertAfter('Artfaktablad. Naturvård – artfakta. ')
$italicRanges += ,@($r.Start, $r.End)
$r = $d.Range($r.End, $r.End)
$r.InsertAfter('SLU Artdatabanken, Uppsala ')

# ---- apply italics recorded above, in a second pass ----
foreach ($rng in $italicRanges) {
    $ir = $d.Range($rng[0], $rng[1])
    $ir.Font.Italic = $true
}

# =====================================================================
# Stage 2: cut the whole scratch block as one range, then re-insert it
# right after "BILAGA 1 - Fridlysta arter" (the last paragraph).
# =====================================================================
$scratchFirstIdx = $scratchStartIdx
$scratchLastIdx = $scratchStartIdx + 12
$blockStart = $d.Paragraphs($scratchFirstIdx).Range.Start
$blockEnd = $d.Paragraphs($scratchLastIdx).Range.End
$block = $d.Range($blockStart, $blockEnd)
$block.Cut()

# BILAGA 1 is Title-styled; a fresh paragraph split off of it would force
# an explicit pStyle=Normal on every pasted paragraph (Paste takes on the
# destination paragraph mark's style). Instead: make one throwaway
# paragraph, explicitly flip IT to Normal (it alone gets an explicit tag,
# and we delete it again below), then split a fresh paragraph off of
# *that* -- which inherits "Normal, no explicit pStyle" cleanly -- and
# paste into it.
$bilaga1 = $d.Paragraphs.Last
$bilaga1.Range.InsertParagraphAfter()
$tempIdx = $d.Paragraphs.Count
$d.Paragraphs($tempIdx).Style = "Normal"
$d.Paragraphs($tempIdx).Range.InsertParagraphAfter()
$pasteIdx = $d.Paragraphs.Count
$d.Paragraphs($pasteIdx).Range.Paste()

# drop the throwaway paragraph
$d.Paragraphs($tempIdx).Range.Delete()

# =====================================================================
# Stage 3: set explicit paragraph styles where the diff calls for them.
# =====================================================================
$d.Paragraphs($tempIdx + 0).Style = "Heading1"
$d.Paragraphs($tempIdx + 6).Style = "Heading2"

# =====================================================================
# Stage 4: header date 2023-09-13 -> 2023-09-15
# The date lives in the "first page" header story (wdHeaderFooterFirstPage
# = 2), which is a separate story range from $d.Content, so Find has to be
# run against that header's own Range.
# =====================================================================
foreach ($sec in $d.Sections) {
    foreach ($hf in $sec.Headers) {
        if ($hf.Exists) {
            $hf.Range.Find.Execute("2023-09-13", $true, $false, $false, $false,
                $false, $true, 1, $false, "2023-09-15", 2) | Out-Null
        }
    }
}
